$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing "Regions"
# header + 76 region names (and their formatting) from column A into
# column B, matching the diff's A-to-B move.
$ws.Columns("A:A").Insert()

# New column A holds the "Language_regions" subset (a sparse list that
# lines up row-for-row with the full region list now living in column B).
# Column B keeps the original list, just renamed from "Regions" to
# "All regions".
$ws.Range("A1").Value2 = "Language_regions"
$ws.Range("B1").Value2 = "All regions"

# Header formatting: A1 bold + centered, B1 keeps the bold style that
# travelled over with the column insert.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter

# Column widths: narrower language-regions column, original width on B.
$ws.Columns("A").ColumnWidth = 30.140625
$ws.Columns("B").ColumnWidth = 45.5703125

$languageRegions = @(
  "Temporal_Mid_R","Temporal_Sup_R","","","","","","",
  "Supp_Motor_Area_R","Frontal_Inf_Oper_R","Frontal_Inf_Tri_R","","",
  "Frontal_Sup_Orb_R","Frontal_Mid_Orb_R","Frontal_Inf_Orb_R","Frontal_Mid_R","",
  "Frontal_Sup_R","Frontal_Sup_L","","Frontal_Mid_L","Frontal_Inf_Orb_L",
  "Frontal_Mid_Orb_L","Frontal_Sup_Orb_L","","","Frontal_Inf_Tri_L",
  "Frontal_Inf_Oper_L","Supp_Motor_Area_L","","","","","","",
  "Temporal_Sup_L","Temporal_Mid_L","Temporal_Inf_L","Precentral_L","",
  "Fusiform_L","","SupraMarginal_L","Angular_L","","Lingual_L","","",
  "Parietal_Inf_L","Parietal_Sup_L","","","Occipital_Inf_L","Occipital_Mid_L",
  "Occipital_Sup_L","Cuneus_L","Cuneus_R","Occipital_Sup_R","Occipital_Mid_R",
  "Occipital_Inf_R","","","Parietal_Sup_R","Parietal_Inf_R","","","Lingual_R",
  "","Angular_R","SupraMarginal_R","","Fusiform_R","","Precentral_R","Temporal_Inf_R"
)

for ($i = 0; $i -lt $languageRegions.Count; $i++) {
    $value = $languageRegions[$i]
    if ($value -ne "") {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value2 = $value
    }
}

# Highlight the two duplicate "first occurrence" rows (Frontal_Sup_L /
# Cuneus_L) with a bold font on a yellow fill.
$highlightRows = @(21, 58)
foreach ($row in $highlightRows) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Bold = $true
    $cell.Interior.Color = 65535   # RGB(255,255,0) / FFFF00
}

# Put the active selection on A4, matching the saved view state.
$ws.Range("A4").Select()
